$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UsuariosRegistro")

$ws.Range("C2").Value = "juan.perez+20251109_022039@test.com"
$ws.Range("C3").Value = "maria.gonzalez+20251109_022039@test.com"
$ws.Range("C4").Value = "carlos.rodriguez+20251109_022039@test.com"
$ws.Range("C5").Value = "ana.martinez+20251109_022039@test.com"
$ws.Range("C6").Value = "luis.garcia+20251109_022039@test.com"
